# SÓC TRĂNG 8-2024 report update
# - new invoice (Cắt mí) data filled into "CHI TIẾT DOANH THU" row 2 + totals
# - "CHI TIẾT VỀ THU NỢ": new "Ngày thực hiện" column inserted before "Lượng thu"
# - "CHI TIẾT CHI TIÊU": three new 08-03-2024 expense entries + updated totals
# - "DOANH SỐ CÁ NHÂN": new sales-person row (Nguyễn Hoàng Yến Quyên) + updated totals
# - "CHI TIÊU TỔNG HỢP": new expense categories + updated totals
# - "LŨY KẾ NGÀY": new day (08-03-2024) rollup row + updated totals
# - "QUỸ LƯƠNG": updated payroll figures
# - "LỢI NHUẬN": reshaped into a per-cơ-sở summary table with new columns

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. CHI TIẾT DOANH THU
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CHI TIẾT DOANH THU")

$ws.Range("E2").Value = "Cắt mí"
$ws.Range("F2").Value = "đường thị út"
$ws.Range("G2").Value = "Cá nhân"
$ws.Range("H2").Value = "Lê Đình Hậu"
$ws.Range("I2").Value = 6000000
$ws.Range("L2").Value = 6000000
$ws.Range("M2").Value = "Nguyễn Hoàng Yến Quyên"
$ws.Range("O2").Value = 6000000
$ws.Range("Q2").Value = 6000000
$ws.Range("S2").Value = "Kha Như Huỳnh "
$ws.Range("U2").Value = 50000
$ws.Range("V2").Value = 0

$ws.Range("I4").Value = 14000000
$ws.Range("L4").Value = 14000000
$ws.Range("O4").Value = 12000000
$ws.Range("Q4").Value = 12000000
$ws.Range("U4").Value = 150000

# ---------------------------------------------------------------------------
# 2. CHI TIẾT VỀ THU NỢ - insert "Ngày thực hiện" column before "Lượng thu"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CHI TIẾT VỀ THU NỢ")
$ws.Range("F1").EntireColumn.Insert()
$ws.Range("F1").Value = "Ngày thực hiện"
$ws.Range("F2").Value = ""

# ---------------------------------------------------------------------------
# 3. CHI TIẾT CHI TIÊU - add 3 new expense rows (CT 752/753/754)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CHI TIẾT CHI TIÊU")
$ws.Range("A3").EntireRow.Insert()
$ws.Range("A3").EntireRow.Insert()
$ws.Range("A3").EntireRow.Insert()

$ws.Range("A3").Value = "CT"
$ws.Range("B3").Value = 752
$ws.Range("C3:C5").NumberFormat = "@"
$ws.Range("C3").Value = "08-03-2024"
$ws.Range("D3").Value = "SÓC TRĂNG"
$ws.Range("E3").Value = "Chi Phí Sinh Hoạt Tại Cơ Sở"
$ws.Range("F3").Value = 180000

$ws.Range("A4").Value = "CT"
$ws.Range("B4").Value = 753
$ws.Range("C4").Value = "08-03-2024"
$ws.Range("D4").Value = "SÓC TRĂNG"
$ws.Range("E4").Value = "Chi Phí Vận Hành"
$ws.Range("F4").Value = 500000

$ws.Range("A5").Value = "CT"
$ws.Range("B5").Value = 754
$ws.Range("C5").Value = "08-03-2024"
$ws.Range("D5").Value = "SÓC TRĂNG"
$ws.Range("E5").Value = "Trang thiết bị Y Tế"
$ws.Range("F5").Value = 700000
$ws.Range("C3:C5").ClearFormats()

$ws.Range("A6").Value = "Tổng"
$ws.Range("B6").Value = 4
$ws.Range("F6").Value = 1560000

# ---------------------------------------------------------------------------
# 4. DOANH SỐ CÁ NHÂN - insert row for Nguyễn Hoàng Yến Quyên (alphabetical slot)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DOANH SỐ CÁ NHÂN")
$ws.Range("A6").EntireRow.Insert()

$ws.Range("A6").Value = "Nguyễn Hoàng Yến Quyên"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 6000000
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0

# row2 (Kha Như Huỳnh) picked up a 2nd phụ phẫu
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 150000

# row5 (Lê Đình Hậu) now sale chính on both invoices
$ws.Range("B5").Value = 14000000

# totals row (now row 11)
$ws.Range("B11").Value = 14000000
$ws.Range("D11").Value = 12000000
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 150000

# ---------------------------------------------------------------------------
# 5. CHI TIÊU TỔNG HỢP - add new expense categories
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CHI TIÊU TỔNG HỢP")
$ws.Range("A4").EntireRow.Insert()
$ws.Range("A4").EntireRow.Insert()

$ws.Range("B2").Value = 360000

$ws.Range("A4").Value = "Trang thiết bị Y Tế"
$ws.Range("B4").Value = 700000

$ws.Range("A3").Value = "Chi Phí Vận Hành"
$ws.Range("B3").Value = 500000

$ws.Range("A5").Value = "Blank"
$ws.Range("B5").Value = 0

$ws.Range("A6").Value = "Tổng cộng"
$ws.Range("B6").Value = 1560000

# ---------------------------------------------------------------------------
# 6. LŨY KẾ NGÀY - add 08-03-2024 rollup
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LŨY KẾ NGÀY")
$ws.Range("A3").EntireRow.Insert()

$ws.Range("B2").Value = 14000000
$ws.Range("C2").Value = 12000000
$ws.Range("G2").Value = 11820000

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "08-03-2024"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1380000
$ws.Range("G3").Value = -1380000

$ws.Range("A4").Value = "Tổng"
$ws.Range("B4").Value = 14000000
$ws.Range("C4").Value = 12000000
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1560000
$ws.Range("G4").Value = 10440000

# ---------------------------------------------------------------------------
# 7. QUỸ LƯƠNG - updated payroll amounts
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("QUỸ LƯƠNG")
$ws.Range("C8").Value = 600000
$ws.Range("C10").Value = 53571.42857142857
$ws.Range("C17").Value = 1380000
$ws.Range("C18").Value = 150000
$ws.Range("C22").Value = 4307380.952380951

# ---------------------------------------------------------------------------
# 8. LỢI NHUẬN - reshape into per-cơ-sở summary table
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LỢI NHUẬN")
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("B1").EntireColumn.Insert()

$ws.Range("A1").Value = "Cơ sở"
$ws.Range("B1").Value = "Tổng đơn giá"
$ws.Range("C1").Value = "Đã thanh toán"
$ws.Range("D1").Value = "Tỉ lệ thanh toán"
$ws.Range("E1").Value = "Tỉ lệ nợ"
$ws.Range("F1").Value = "Thu nợ"
$ws.Range("G1").Value = "Tổng doanh thu"
$ws.Range("H1").Value = "Chi tiêu"
$ws.Range("I1").Value = "Quỹ lương"
$ws.Range("J1").Value = "Tổng chi phí"
$ws.Range("K1").Value = "Lợi nhuận"
$ws.Range("L1").Value = "Tỉ lệ lợi nhuận"

$ws.Range("A2").Value = "SÓC TRĂNG"
$ws.Range("B2").Value = 7000000
$ws.Range("C2").Value = 6000000
$ws.Range("D2").Value = 0.8571428571428571
$ws.Range("E2").Value = 0.1428571428571429
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 6000000
$ws.Range("H2").Value = 780000
$ws.Range("I2").Value = 4307380.952380951
$ws.Range("J2").Value = 5087380.952380951
$ws.Range("K2").Value = 912619.0476190485
$ws.Range("L2").Value = 0.1521031746031747

Write-Host "edit applied"
